# Improve numerical stability on input data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the S2 column (column D) entirely; its values move into column C.
$ws.Columns.Item(4).Delete()

# Update the probability values for the existing rows (years 2023-2040).
# Column B = S0, Column C = S1 (now holding what used to be the S2 values).
$ws.Range("B2").Value = 0.984
$ws.Range("C2").Value = 0.016

$ws.Range("B3").Value = 0.96
$ws.Range("C3").Value = 0.04

$ws.Range("B4").Value = 0.96
$ws.Range("C4").Value = 0.04

$ws.Range("B5").Value = 0.96
$ws.Range("C5").Value = 0.04

$ws.Range("B6").Value = 0.96
$ws.Range("C6").Value = 0.04

# Append new rows for years 2045 and 2050 with the same probabilities.
$ws.Range("A7").Value = 2045
$ws.Range("B7").Value = 0.96
$ws.Range("C7").Value = 0.04

$ws.Range("A8").Value = 2050
$ws.Range("B8").Value = 0.96
$ws.Range("C8").Value = 0.04
